$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner text (A1)
$ws.Range("A1").Value = "Datos actualizados a 10 de Abril de 2020 a las 01:52"

# Row 4 - Estados Unidos (country unchanged, stats updated)
$ws.Range("B4").Value = 466969
$ws.Range("C4").Value = 31939
$ws.Range("D4").Value = 25316
$ws.Range("E4").Value = 425017
$ws.Range("F4").Value = 10011
$ws.Range("G4").Value = 1845
$ws.Range("H4").Value = 16636

# Row 17 - Brasil (country unchanged, stats updated)
$ws.Range("B17").Value = 18145
$ws.Range("C17").Value = 1957
$ws.Range("E17").Value = 17018
$ws.Range("G17").Value = 134
$ws.Range("H17").Value = 954

# Row 26 - Noruega (country unchanged, stats updated)
$ws.Range("B26").Value = 6219
$ws.Range("C26").Value = 177
$ws.Range("E26").Value = 6079

# Row 27 - Australia (country unchanged, stats updated)
$ws.Range("B27").Value = 6152
$ws.Range("C27").Value = 100
$ws.Range("E27").Value = 3114

# Panama overtakes Finlandia in the ranking: Panama's updated stats move into
# row 45, Finlandia's (unchanged) stats shift down into row 46.
$ws.Range("A45").Value = "Panama"
$ws.Range("B45").Value = 2752
$ws.Range("C45").Value = 224
$ws.Range("D45").Value = 16
$ws.Range("E45").Value = 2670
$ws.Range("F45").Value = 107
$ws.Range("G45").Value = 3
$ws.Range("H45").Value = 66

$ws.Range("A46").Value = "Finlandia"
$ws.Range("B46").Value = 2605
$ws.Range("C46").Value = 118
$ws.Range("D46").Value = 300
$ws.Range("E46").Value = 2263
$ws.Range("F46").Value = 82
$ws.Range("G46").Value = 2
$ws.Range("H46").Value = 42

# Argentina overtakes Ucrania in the ranking: Argentina's updated stats move
# into row 54, Ucrania's (unchanged) stats shift down into row 55.
$ws.Range("A54").Value = "Argentina"
$ws.Range("B54").Value = 1894
$ws.Range("C54").Value = 99
$ws.Range("D54").Value = 365
$ws.Range("E54").Value = 1450
$ws.Range("F54").Value = 96
$ws.Range("G54").Value = 14
$ws.Range("H54").Value = 79

$ws.Range("A55").Value = "Ucrania"
$ws.Range("B55").Value = 1892
$ws.Range("C55").Value = 224
$ws.Range("D55").Value = 45
$ws.Range("E55").Value = 1790
$ws.Range("F55").Value = 33
$ws.Range("G55").Value = 5
$ws.Range("H55").Value = 57

# Row 123 - Republica de Yibuti (country unchanged, stats updated)
$ws.Range("B123").Value = 140
$ws.Range("C123").Value = 5
$ws.Range("D123").Value = 28
$ws.Range("E123").Value = 111

# Row 128 - Ruanda (country unchanged, stats updated)
$ws.Range("B128").Value = 113
$ws.Range("C128").Value = 3
$ws.Range("E128").Value = 106
